$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.990.65'
$ws.Range('E2').Value = '  +2.89%  '

$ws.Range('D3').Value = '2.050.76'
$ws.Range('E3').Value = '  +2.32%  '

$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '229.41'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.96%  '

$ws.Range('E6').Value = '  +2.62%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '58.43'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +7.15%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('E9').Value = '  +3.18%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0809'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.35%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.103'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.16%  '

$ws.Range('D12').Value = '2.354.37'
$ws.Range('E12').Value = '  +2.25%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '14.57'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +4.51%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '20.79'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +5.47%  '

$ws.Range('E15').Value = '  +2.46%  '

$ws.Range('E16').Value = '  +0.75%  '

$ws.Range('D17').Value = '2.068.66'
$ws.Range('E17').Value = '  +0.86%  '

$ws.Range('D18').Value = '37.896.26'
$ws.Range('E18').Value = '  +2.88%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.35'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.78%  '

$ws.Range('E20').Value = '  +2.10%  '

$ws.Range('D21').Value = '0.0₃0835'
$ws.Range('E21').Value = '  +3.05%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '224.29'
$ws.Range('D22').ClearFormats()

$ws.Range('E23').Value = '  +0.06%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.43'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.04%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.24'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +4.23%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '166.49'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.59%  '

$ws.Range('E27').Value = '  +2.93%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.132'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +5.98%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.97'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.66%  '

$ws.Range('E30').Value = '  +2.46%  '

$ws.Range('E31').Value = '  +2.69%  '

$ws.Range('E32').Value = '  +1.88%  '

$ws.Range('E33').Value = '  +10.88%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.57'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.60%  '

$ws.Range('E35').Value = '  +1.43%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.31'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.20%  '

$ws.Range('E37').Value = '  +14.03%  '

$ws.Range('E38').Value = '  +5.70%  '

$ws.Range('E39').Value = '  +0.02%  '

$ws.Range('D40').Value = '1.519.74'
$ws.Range('E40').Value = '  +4.66%  '

$ws.Range('B41').Value = 'HuobiToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.89'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.91%  '

$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0217'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.88%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '96.89'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.68%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.50'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +4.19%  '

$ws.Range('E45').Value = '  +1.25%  '

$ws.Range('E46').Value = '  +0.94%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.08'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +16.16%  '

$ws.Range('E48').Value = '  +2.26%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.96'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.48%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.09'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.25%  '

$ws.Range('D51').Value = '2.240.75'
$ws.Range('E51').Value = '  +2.18%  '
